$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A (read-only key/value view needs more room for labels) ---
$ws.Columns("A").ColumnWidth = 33.0

# --- Row 1: extend header with KV2..KV10 ---
$ws.Range("Q1").Value = "KV2"
$ws.Range("R1").Value = "KV3"
$ws.Range("S1").Value = "KV4"
$ws.Range("T1").Value = "KV5"
$ws.Range("U1").Value = "KV6"
$ws.Range("V1").Value = "KV7"
$ws.Range("W1").Value = "KV8"
$ws.Range("X1").Value = "KV9"
$ws.Range("Y1").Value = "KV10"

# --- Row 2: second key/value pair ---
$ws.Range("Q2").Value = "two"

# --- Row 3: second key/value pair (numeric) ---
$ws.Range("Q3").Value = 456

# --- Row 4: extra key/value text, formatted like the date cell next to it ---
$ws.Range("Q4").NumberFormat = "DD/MM/YY"
$ws.Range("Q4").Value = "text text text"

# --- Row 5: new record - "TEN KVs test" with all ten key/value pairs filled ---
$ws.Range("A5").Value = "2018-08-11T17:26:56.812+06"
$ws.Range("B5").Value = "2018-08-11T17:29:09.865+06"
$ws.Range("C5").Value = "viviane_lucia"
$ws.Range("D5").Value = "TEN KVs test"
$ws.Range("E5").Value = "2018-08-11T17:27:00.000+06"
$ws.Range("F5").Value = "female"
$ws.Range("G5").Value = "Age 19-25"
$ws.Range("H5").Value = "Kutupalong RC"
$ws.Range("I5").Value = "English"
$ws.Range("J5").Value = 3
$ws.Range("K5").Font.Bold = $true
$ws.Range("K5").Font.Bold = $false
$ws.Range("L5").Value = 27671657
$ws.Range("M5").Value = "553c960d-3c2f-4349-a843-742cc9495dc8"
$ws.Range("N5").Value = "2018-08-11T11:29:26"
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = "testKV1"
$ws.Range("Q5").Value = "testKV2"
$ws.Range("R5").Value = "testKV3"
$ws.Range("S5").Value = "testKV4"
$ws.Range("T5").Value = "testKV5"
$ws.Range("U5").Value = "testKV6"
$ws.Range("V5").Value = "testKV7"
$ws.Range("W5").Value = "testKV8"
$ws.Range("X5").Value = "testKV9"
$ws.Range("Y5").Value = "testKV10"

# --- Row 6: new record - "NO Kvs test" with no key/value pairs ---
$ws.Range("A6").Value = "2018-08-11T17:26:57.812+06"
$ws.Range("B6").Value = "2018-08-11T17:29:09.865+06"
$ws.Range("C6").Value = "viviane_lucia"
$ws.Range("D6").Value = "NO Kvs test"
$ws.Range("E6").Value = "2018-08-11T17:27:00.000+06"
$ws.Range("F6").Value = "female"
$ws.Range("G6").Value = "Age 19-25"
$ws.Range("H6").Value = "Kutupalong RC"
$ws.Range("I6").Value = "English"
$ws.Range("J6").Value = 3
$ws.Range("K6").Font.Bold = $true
$ws.Range("K6").Font.Bold = $false
$ws.Range("L6").Value = 27671657
$ws.Range("M6").Value = "553c960d-3c2f-4349-a843-742cc9495dc8"
$ws.Range("N6").Value = "2018-08-11T11:29:26"
$ws.Range("O6").Value = 3

# --- Leave the view on the last-entered row, as in the read-only item view ---
$ws.Range("D6").Select()
